$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quotation row appended at the bottom of the table (row 54).
# Column A holds a date-looking string ("2025-10-02"); Excel's COM layer
# auto-coerces such text into a date serial on assignment, so we force a
# Text number format just long enough to land the literal string, then
# clear the formatting again so the cell's style stays the default (style 0),
# matching the rest of the sheet (which carries no explicit styles).
$cellA = $ws.Range("A54")
$cellA.NumberFormat = "@"
$cellA.Value = "2025-10-02"
$cellA.ClearFormats()

$ws.Range("B54").Value = "21:20:54"
$ws.Range("C54").Value = "1.00 EUR = 1,834.9959"
